# Update '想去人数' (want-to-go count) values in column F across sheets
# as refreshed from the bilibili source data (per commit: gh-pages output regenerated).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 85   # was 83
$ws.Range("F9").Value = 2600   # was 2598
$ws.Range("F15").Value = 337   # was 336
$ws.Range("F18").Value = 2130   # was 2129
$ws.Range("F22").Value = 2624   # was 2622
$ws.Range("F23").Value = 4   # was 3
$ws.Range("F28").Value = 421   # was 419
$ws.Range("F37").Value = 4587   # was 4586
$ws.Range("F38").Value = 155   # was 154

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 4204   # was 4202
$ws.Range("F13").Value = 317   # was 316
$ws.Range("F14").Value = 327   # was 326
$ws.Range("F18").Value = 155   # was 154
$ws.Range("F29").Value = 271   # was 270
$ws.Range("F38").Value = 479   # was 478

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F6").Value = 532   # was 531
$ws.Range("F7").Value = 133   # was 132

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 532   # was 531
$ws.Range("F5").Value = 85   # was 83
$ws.Range("F12").Value = 133   # was 132
$ws.Range("F16").Value = 2600   # was 2598
$ws.Range("F22").Value = 337   # was 336
$ws.Range("F25").Value = 327   # was 326
$ws.Range("F27").Value = 2130   # was 2129
$ws.Range("F30").Value = 155   # was 154
$ws.Range("F32").Value = 2624   # was 2622
$ws.Range("F40").Value = 421   # was 419
$ws.Range("F41").Value = 421   # was 419
$ws.Range("F48").Value = 4587   # was 4586
$ws.Range("F49").Value = 155   # was 154
$ws.Range("F50").Value = 479   # was 478
